# Swap the contents of column C (codeforiati:group-name) and column D
# (codeforiati:group-code), including the header row, so that column C
# becomes "codeforiati:group-code" and column D becomes
# "codeforiati:group-name" for every row of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)

    $cValue = $cCell.Value2
    $dValue = $dCell.Value2

    $cCell.Value2 = $dValue
    $dCell.Value2 = $cValue
}
